# [Improvement] On terminology : room -> bed

$wb = $excel.ActiveWorkbook

# The "rooms" sheet is being renamed to "beds"
$roomsSheet = $wb.Worksheets.Item("rooms")

# Update the header row wording from *rooms* to *beds* before renaming the sheet
$roomsSheet.Range("A1").Value = "all_beds"
$roomsSheet.Range("B1").Value = "new_beds"
$roomsSheet.Range("C1").Value = "old_beds"
$roomsSheet.Range("E1").Value = "new_beds_service"
$roomsSheet.Range("F1").Value = "old_beds_service"
$roomsSheet.Range("G1").Value = "beds_capacities"

# Rename the sheet itself
$roomsSheet.Name = "beds"

# Move the selection on the "beds" sheet to F6 and make it the active sheet/tab
$roomsSheet.Activate()
$roomsSheet.Range("F6").Select()
